# LOQ4088.xlsx edit: fill in real syllabus content for the "Docentes
# responsáveis", "Programa resumido", "Programa", "Método", "Critério",
# "Norma de recuperação" and "Bibliografia" rows, plus update "Objetivos:".
#
# Before the edit, rows 13-23 held placeholder/duplicated text (several
# rows literally repeated values copy-pasted from other rows, e.g. row 13
# "Programa resumido" showed "Semestral", row 18 "Método" showed the
# professor's name, etc.). The real content is inserted here, which also
# requires one brand-new row (a dedicated "Docentes responsáveis:" value
# row), pushing every row from 13 on down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new row used for "Docentes responsáveis:" value ----------
# (this shifts old rows 13..23 down to 14..24, carrying their row heights
# and A-column labels with them)
$ws.Rows(13).Insert()

# Row 13 used to carry the A-column "Programa resumido:" label; now that
# belongs on row 14, so row 13's A cell must go away, leaving only B/C.
$ws.Range("A13").Clear()

# --- row 10: Objetivos: ---------------------------------------------------
$objetivos = "Aplicar os conceitos fundamentais relacionados aos processos físicos químicos, ampliando o conhecimento termodinâmico dos sistemas, isto é, a definição dos critérios de equilíbrio e de espontaneidade para misturas e reações químicas."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- row 13: Docentes responsáveis: (value only, new row) ----------------
$docentes = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes
# Match the sheet's normal "value cell" formatting (wrap text, regular
# weight, top aligned) instead of the bold label look the row inherited
# from the insert.
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

# --- row 14: Programa resumido: -------------------------------------------
$programaResumido = "Termodinâmica de soluções. Equilíbrio líquido  vapor. Equilíbrio de fases. Equilíbrio em reações químicas  Equilíbrio químico"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- row 16: Programa: -----------------------------------------------------
$programa = @"
1- Termodinâmica de soluções 
1.1- Relações fundamentais entre propriedades 
1.2- O potencial químico 
1.3- Fugacidade e coeficiente de fugacidade 
1.4- A solução Ideal 
1.5- Modelos para a energia de Gibbs 
1.6- Propriedades de mistura 
1.7- Efeitos térmicos em processos de mistura 
2- Equilíbrio liquido  vapor 
2.1- A natureza em equilíbrio 
2.2- A regra das fases. Teorema de Duhem 
2.3- Calculo dos pontos de orvalho e de bolha 
2.4- Calculo de Flash 
3- Equilíbrio de fases 
3.1- Equilíbrio e estabilidade 
3.2- Equilíbrio líquido-líquido 
3.3- Equilíbrio líquido-líquido-vapor 
3.4- Equilíbrio sólido-líquido 
3.5- Equilíbrio sólido-vapor 
3.6- Equilíbrio na adsorção de gases em sólidos 
4- Equilíbrio em reações químicas  Equilíbrio químico 
4.1- A variação de energia de Gibbs padrão e a constante de equilíbrio 
4.2- Efeito da temperatura sobre a constante de equilíbrio 
4.3- Avaliação das constantes de equilíbrio 
4.4- Relação entre as constantes de equilíbrio e a composição 
4.5- Conversões de equilíbrio em reações isoladas
"@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa
$ws.Rows(16).RowHeight = 120

# --- row 19: Método: ---------------------------------------------------
$metodo = "A avaliação será feita por meio de duas provas escritas (P1 e P2)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- row 20: Critério: --------------------------------------------------
$criterio = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- row 21: Norma de recuperação: ---------------------------------------
$norma = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- row 22: Bibliografia: ------------------------------------------------
$biblio = @"
KORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007. 
MORAN, M. J.; SHAPIRO, H. N. Princípios de Termodinâmica para Engenharia. 1 ed. LTC Editora, 2009. 
SANDLER, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999 
SMITH, J.M.; VAN NESS, H.C.; Abott, M. M. Introdução à Termodinâmica da Engenharia Química. 7ª ed. LTC editora, 2007. 
TERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009. 
VAN WILEN, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica. 6 ed. 2004
"@
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
